$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 385, shifting the existing rows 385-399 down to 386-400
# (new weekly price observation added at the top of the Choclo data block).
$ws.Rows.Item(385).Insert()

$ws.Cells.Item(385, 1).Value = 3
$ws.Cells.Item(385, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(385, 3).Value = "Coquimbo"
$ws.Cells.Item(385, 4).Value = 44509
$ws.Cells.Item(385, 5).Value = 5
$ws.Cells.Item(385, 6).Value = 100112024
$ws.Cells.Item(385, 7).Value = "Choclo"
$ws.Cells.Item(385, 8).Value = "Dulce o Americano"
$ws.Cells.Item(385, 9).Value = "Primera"
$ws.Cells.Item(385, 10).Value = 160
$ws.Cells.Item(385, 11).Value = 40000
$ws.Cells.Item(385, 12).Value = 41000
$ws.Cells.Item(385, 13).Value = 40656
$ws.Cells.Item(385, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(385, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(385, 16).Value = 581
$ws.Cells.Item(385, 17).Value = 70
$ws.Cells.Item(385, 18).Value = "Hortaliza"
